$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2032693333333333
$ws.Range("H2").Value = 0.609808
$ws.Range("I2").Value = 0.01014574830092973
$ws.Range("J2").Value = 0.01014574830092973
$ws.Range("M2").Value = 0.4070896666666666
$ws.Range("N2").Value = 1.221269
$ws.Range("O2").Value = 0.06732963386756073
$ws.Range("P2").Value = 0.06732963386756075
$ws.Range("Q2").Value = 0.08274884515022221
$ws.Range("R2").Value = 0.744739606352
$ws.Range("S2").Value = 0.000683109518414025
$ws.Range("T2").Value = 0.0006831095184140251

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2032693333333333
$ws.Range("H3").Value = 0.609808
$ws.Range("I3").Value = 0.01014574830092973
$ws.Range("J3").Value = 0.01014574830092973
$ws.Range("M3").Value = 2.455660333333333
$ws.Range("N3").Value = 7.366980999999999
$ws.Range("O3").Value = 0.4061481405319192
$ws.Range("P3").Value = 0.4061481405319192
$ws.Range("Q3").Value = 0.4991604388497777
$ws.Range("R3").Value = 4.492443949648
$ws.Range("S3").Value = 0.004120676806727487
$ws.Range("T3").Value = 0.004120676806727488

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2032693333333333
$ws.Range("H4").Value = 0.609808
$ws.Range("I4").Value = 0.01014574830092973
$ws.Range("J4").Value = 0.01014574830092973
$ws.Range("M4").Value = 3.183468333333334
$ws.Range("N4").Value = 9.550405000000001
$ws.Range("O4").Value = 0.52652222560052
$ws.Range("P4").Value = 0.5265222256005201
$ws.Range("Q4").Value = 0.6471014858044445
$ws.Range("R4").Value = 5.823913372240001
$ws.Range("S4").Value = 0.005341961975788215
$ws.Range("T4").Value = 0.005341961975788216

# Row 5
$ws.Range("I5").Value = 0.9454503978493691
$ws.Range("J5").Value = 0.9454503978493691
$ws.Range("M5").Value = 0.4070896666666666
$ws.Range("N5").Value = 1.221269
$ws.Range("O5").Value = 0.06732963386756073
$ws.Range("P5").Value = 0.06732963386756075
$ws.Range("Q5").Value = 7.711104814386554
$ws.Range("R5").Value = 69.399943329479
$ws.Range("S5").Value = 0.06365682912713765
$ws.Range("T5").Value = 0.06365682912713766

# Row 6
$ws.Range("I6").Value = 0.9454503978493691
$ws.Range("J6").Value = 0.9454503978493691
$ws.Range("M6").Value = 2.455660333333333
$ws.Range("N6").Value = 7.366980999999999
$ws.Range("O6").Value = 0.4061481405319192
$ws.Range("P6").Value = 0.4061481405319192
$ws.Range("Q6").Value = 46.51519252236344
$ws.Range("R6").Value = 418.6367327012709
$ws.Range("S6").Value = 0.3839929210516845
$ws.Range("T6").Value = 0.3839929210516845

# Row 7
$ws.Range("I7").Value = 0.9454503978493691
$ws.Range("J7").Value = 0.9454503978493691
$ws.Range("M7").Value = 3.183468333333334
$ws.Range("N7").Value = 9.550405000000001
$ws.Range("O7").Value = 0.52652222560052
$ws.Range("P7").Value = 0.5265222256005201
$ws.Range("Q7").Value = 60.30135373520611
$ws.Range("R7").Value = 542.7121836168551
$ws.Range("S7").Value = 0.4978006476705469
$ws.Range("T7").Value = 0.497800647670547

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.889628
$ws.Range("H8").Value = 2.668884
$ws.Range("I8").Value = 0.04440385384970112
$ws.Range("J8").Value = 0.04440385384970111
$ws.Range("M8").Value = 0.4070896666666666
$ws.Range("N8").Value = 1.221269
$ws.Range("O8").Value = 0.06732963386756073
$ws.Range("P8").Value = 0.06732963386756075
$ws.Range("Q8").Value = 0.3621583659773333
$ws.Range("R8").Value = 3.259425293796
$ws.Range("S8").Value = 0.002989695222009053
$ws.Range("T8").Value = 0.002989695222009053

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.889628
$ws.Range("H9").Value = 2.668884
$ws.Range("I9").Value = 0.04440385384970112
$ws.Range("J9").Value = 0.04440385384970111
$ws.Range("M9").Value = 2.455660333333333
$ws.Range("N9").Value = 7.366980999999999
$ws.Range("O9").Value = 0.4061481405319192
$ws.Range("P9").Value = 0.4061481405319192
$ws.Range("Q9").Value = 2.184624191022666
$ws.Range("R9").Value = 19.661617719204
$ws.Range("S9").Value = 0.01803454267350721
$ws.Range("T9").Value = 0.01803454267350721

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.889628
$ws.Range("H10").Value = 2.668884
$ws.Range("I10").Value = 0.04440385384970112
$ws.Range("J10").Value = 0.04440385384970111
$ws.Range("M10").Value = 3.183468333333334
$ws.Range("N10").Value = 9.550405000000001
$ws.Range("O10").Value = 0.52652222560052
$ws.Range("P10").Value = 0.5265222256005201
$ws.Range("Q10").Value = 2.832102566446667
$ws.Range("R10").Value = 25.48892309802
$ws.Range("S10").Value = 0.02337961595418485
$ws.Range("T10").Value = 0.02337961595418485

